# Edit date in godisnji.docx
# "Đurđevac, 3.6.2024." -> "Đurđevac, 7. 6. 2024."
#
# The target OOXML keeps the leading "Đurđevac, " text in the original
# run, and represents the new date as a sequence of separate runs
# (same run properties) for "7", ".", " ", "6.", " ", "2024." -- as if
# each piece had been edited/typed individually. We first perform the
# textual replacement (as a single Find/Replace, which keeps everything
# in one run), and only afterwards split that run into the matching
# pieces -- touching each piece's formatting (no-op Bold toggle) forces
# the engine to materialise it as its own <w:r> without altering the
# text. Splitting must happen last: any further text edit after a split
# re-coalesces adjacent identically-formatted runs back together.

$d = $word.ActiveDocument

$d.Content.Find.Execute("Đurđevac, 3.6.2024.", $true, $false, $false, $false, $false, $true, 1, $false, "Đurđevac, 7. 6. 2024.", 2)

$r = $d.Content
$r.Find.Execute("Đurđevac, 7. 6. 2024.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $r.Start

$pieces = @("Đurđevac, ", "7", ".", " ", "6.", " ", "2024.")
$pos = $matchStart
foreach ($piece in $pieces) {
    $pieceStart = $pos
    $pieceEnd = $pos + $piece.Length
    $seg = $d.Range($pieceStart, $pieceEnd)
    $seg.Font.Bold = $true
    $seg.Font.Bold = $false
    $pos = $pieceEnd
}
